$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column C ("Severity"), shifting the old
#    "Issues Raised" (C->D) and "Resolution Status" (D->E) columns right.
$ws.Columns("C").Insert()

# The inserted column doesn't inherit column B's explicit width, so match it
# up (column C ends up the same width as B, per the saved workbook).
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth()

# 2. Populate the new header + first data row for the Severity column, and
#    update the two cells whose displayed text changed.
$ws.Range("B3").Value = "Passed"
$ws.Range("E2").Value = "Open"
$ws.Range("C1").Value = "Severity"
$ws.Range("C2").Value = "Low"

# 3. Rebuild the data validation rules.
#    Before the column insert there were two rules:
#      B2:B12 -> "Pass, Failed"
#      D2:D12 -> " , In Progress, Resolved, "   (now E2:E12 after the shift)
#    Clear those two rules completely, then lay down the final set of rules
#    fresh (mirrors how the rules ended up split across the sheet once the
#    "Severity" and re-worked "Passed"/"Open" values were introduced).
$ws.Range("B2:C12").Validation.Delete()
$ws.Range("E2:E12").Validation.Delete()

$ws.Range("B4:C12").Validation.Add(3, 1, 1, '"Pass, Failed"')
$ws.Range("E3:E12").Validation.Add(3, 1, 1, '" , In Progress, Resolved, "')
$ws.Range("B3:C3").Validation.Add(3, 1, 1, '"Passed, Failed"')
$ws.Range("B2").Validation.Add(3, 1, 1, '"Passed, Failed"')
$ws.Range("E2").Validation.Add(3, 1, 1, '"Open, In Progress, Resolved, Reopened"')
$ws.Range("C2").Validation.Add(3, 1, 1, '"Low, Medium, High"')

# 4. Match the saved selection in the source workbook.
$ws.Range("C3").Select()
